$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Execution Required" column (F) for rows 4 through 21 from "Yes" to "No"
for ($row = 4; $row -le 21; $row++) {
    $ws.Cells.Item($row, 6).Value = "No"
}

# Update the selection to span F3:F21 (keeping F3 as the active cell)
$ws.Range("F3:F21").Select()
